$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - update the "reps" header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) - update values and delete C2
$ws.Range("B2").Value = 398.61166425469929
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 393.25617965248244
$ws.Range("E2").Value = 626.63208474871317

# Row 3 (STR) - delete B3, update C3, add D3, update E3
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 479.48374785683899
$ws.Range("D3").Value = 402.15337977147959
$ws.Range("E3").Value = 416.48106036082805

# Update selection to match new range
$ws.Range("B1:E3").Select()
